# Auto-generated edit script: update cached market-data snapshot values
# across the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 2482.1538  # H98: was 2585.45
$ws.Cells.Item(98, 9).Value = 2606.4119  # I98: was 2720.9143
$ws.Cells.Item(98, 11).Value = 2606.4119  # K98: was 2720.9143
$ws.Cells.Item(98, 13).Value = -1108.4119  # M98: was -1222.9143
$ws.Cells.Item(112, 8).Value = 1842.7084  # H112: was 1913.1555
$ws.Cells.Item(112, 10).Value = 2001.2195  # J112: was 2097.158
$ws.Cells.Item(112, 12).Value = 6003.6585  # L112: was 6291.474
$ws.Cells.Item(112, 14).Value = -8219.6585  # N112: was -8507.474
$ws.Cells.Item(122, 8).Value = 2482.1538  # H122: was 2585.45
$ws.Cells.Item(122, 9).Value = 2606.4119  # I122: was 2720.9143
$ws.Cells.Item(122, 11).Value = 7819.2357  # K122: was 8162.742899999999
$ws.Cells.Item(122, 13).Value = -5369.2357  # M122: was -5712.742899999999
$ws.Cells.Item(132, 8).Value = 7412990.5  # H132: was 6541184.5
$ws.Cells.Item(132, 9).Value = 8776339  # I132: was 9263910
$ws.Cells.Item(132, 10).Value = 11953.429  # J132: was 6644.933
$ws.Cells.Item(132, 11).Value = 26329017  # K132: was 27791730
$ws.Cells.Item(132, 12).Value = 35860.287  # L132: was 19934.799
$ws.Cells.Item(132, 13).Value = -26326487  # M132: was -27789200
$ws.Cells.Item(132, 14).Value = -40920.287  # N132: was -24994.799
$ws.Cells.Item(137, 8).Value = 1291.6595  # H137: was 1052.9565
$ws.Cells.Item(137, 9).Value = 887.5599999999999  # I137: was 714.2195
$ws.Cells.Item(137, 10).Value = 1750.8636  # J137: was 1548.9642
$ws.Cells.Item(137, 11).Value = 2662.68  # K137: was 2142.6585
$ws.Cells.Item(137, 12).Value = 5252.5908  # L137: was 4646.892599999999
$ws.Cells.Item(137, 13).Value = -112.6799999999998  # M137: was 407.3415
$ws.Cells.Item(137, 14).Value = -10352.5908  # N137: was -9746.892599999999
$ws.Cells.Item(138, 8).Value = 1331.6632  # H138: was 1351.7576
$ws.Cells.Item(138, 9).Value = 884.2632  # I138: was 901.4054
$ws.Cells.Item(138, 10).Value = 1629.9298  # J138: was 1620.5161
$ws.Cells.Item(138, 11).Value = 2652.7896  # K138: was 2704.2162
$ws.Cells.Item(138, 12).Value = 4889.7894  # L138: was 4861.5483
$ws.Cells.Item(138, 13).Value = 2487.2104  # M138: was 2435.7838
$ws.Cells.Item(138, 14).Value = -15169.7894  # N138: was -15141.5483
$ws.Cells.Item(141, 8).Value = 622.04  # H141: was 435.6389
$ws.Cells.Item(141, 9).Value = 564.625  # I141: was 435.6389
$ws.Cells.Item(141, 10).Value = 2000  # J141: was 0
$ws.Cells.Item(141, 11).Value = 1693.875  # K141: was 1306.9167
$ws.Cells.Item(141, 12).Value = 6000  # L141: was 0
$ws.Cells.Item(141, 13).Value = 3486.125  # M141: was 3873.0833
$ws.Cells.Item(141, 14).Value = -16360  # N141: was None

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4016.9614  # H32: was 3973.175
$ws.Cells.Item(32, 9).Value = 3525.4429  # I32: was 3640.9412
$ws.Cells.Item(32, 10).Value = 8317.75  # J32: was 5855.8335
$ws.Cells.Item(32, 11).Value = 3525.4429  # K32: was 3640.9412
$ws.Cells.Item(32, 12).Value = 8317.75  # L32: was 5855.8335
$ws.Cells.Item(32, 13).Value = -3238.4429  # M32: was -3353.9412
$ws.Cells.Item(32, 14).Value = -8891.75  # N32: was -6429.8335
$ws.Cells.Item(61, 8).Value = 18519502  # H61: was 17544796
$ws.Cells.Item(61, 9).Value = 25000646  # I61: was 24390864
$ws.Cells.Item(61, 10).Value = 1950.8572  # J61: was 1747
$ws.Cells.Item(61, 11).Value = 25000646  # K61: was 24390864
$ws.Cells.Item(61, 12).Value = 1950.8572  # L61: was 1747
$ws.Cells.Item(61, 13).Value = -25000434  # M61: was -24390652
$ws.Cells.Item(61, 14).Value = -2374.8572  # N61: was -2171
$ws.Cells.Item(74, 8).Value = 1150.4117  # H74: was 1256.2903
$ws.Cells.Item(74, 9).Value = 815.84  # I74: was 931.6087
$ws.Cells.Item(74, 10).Value = 2079.7778  # J74: was 2189.75
$ws.Cells.Item(74, 11).Value = 815.84  # K74: was 931.6087
$ws.Cells.Item(74, 12).Value = 2079.7778  # L74: was 2189.75
$ws.Cells.Item(74, 13).Value = 58.15999999999997  # M74: was -57.6087
$ws.Cells.Item(74, 14).Value = -3827.7778  # N74: was -3937.75
$ws.Cells.Item(77, 8).Value = 1150.4117  # H77: was 1256.2903
$ws.Cells.Item(77, 9).Value = 815.84  # I77: was 931.6087
$ws.Cells.Item(77, 10).Value = 2079.7778  # J77: was 2189.75
$ws.Cells.Item(77, 11).Value = 4079.2  # K77: was 4658.0435
$ws.Cells.Item(77, 12).Value = 10398.889  # L77: was 10948.75
$ws.Cells.Item(77, 13).Value = 288.7999999999997  # M77: was -290.0434999999998
$ws.Cells.Item(77, 14).Value = -19134.889  # N77: was -19684.75
$ws.Cells.Item(112, 8).Value = 10171.8  # H112: was 10346.5
$ws.Cells.Item(112, 10).Value = 10171.8  # J112: was 10346.5
$ws.Cells.Item(112, 12).Value = 10171.8  # L112: was 10346.5
$ws.Cells.Item(112, 14).Value = -13125.8  # N112: was -13300.5
$ws.Cells.Item(122, 8).Value = 1995.625  # H122: was 1725.5
$ws.Cells.Item(122, 10).Value = 1440  # J122: was 910
$ws.Cells.Item(122, 12).Value = 4320  # L122: was 2730
$ws.Cells.Item(122, 14).Value = -9220  # N122: was -7630
$ws.Cells.Item(132, 8).Value = 1076.1356  # H132: was 1333.8909
$ws.Cells.Item(132, 9).Value = 1009.1395  # I132: was 1128.8379
$ws.Cells.Item(132, 10).Value = 1256.1875  # J132: was 1755.3889
$ws.Cells.Item(132, 11).Value = 3027.4185  # K132: was 3386.5137
$ws.Cells.Item(132, 12).Value = 3768.5625  # L132: was 5266.1667
$ws.Cells.Item(132, 13).Value = -497.4184999999998  # M132: was -856.5137
$ws.Cells.Item(132, 14).Value = -8828.5625  # N132: was -10326.1667
$ws.Cells.Item(136, 8).Value = 18519502  # H136: was 17544796
$ws.Cells.Item(136, 9).Value = 25000646  # I136: was 24390864
$ws.Cells.Item(136, 10).Value = 1950.8572  # J136: was 1747
$ws.Cells.Item(136, 11).Value = 75001938  # K136: was 73172592
$ws.Cells.Item(136, 12).Value = 5852.571599999999  # L136: was 5241
$ws.Cells.Item(136, 13).Value = -74999388  # M136: was -73170042
$ws.Cells.Item(136, 14).Value = -10952.5716  # N136: was -10341

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(51, 8).Value = 0  # H51: was 45000
$ws.Cells.Item(51, 10).Value = 0  # J51: was 45000
$ws.Cells.Item(51, 12).Value = 0  # L51: was 45000
$ws.Cells.Item(51, 14).ClearContents()  # N51: remove (was -45982)
$ws.Cells.Item(57, 8).Value = 26691.334  # H57: was 25220
$ws.Cells.Item(57, 9).Value = 0  # I57: was 25000
$ws.Cells.Item(57, 10).Value = 26691.334  # J57: was 25330
$ws.Cells.Item(57, 11).Value = 0  # K57: was 25000
$ws.Cells.Item(57, 12).Value = 26691.334  # L57: was 25330
$ws.Cells.Item(57, 13).ClearContents()  # M57: remove (was -24280)
$ws.Cells.Item(57, 14).Value = -28131.334  # N57: was -26770
$ws.Cells.Item(94, 8).Value = 11364303  # H94: was 8333852
$ws.Cells.Item(94, 9).Value = 20833616  # I94: was 13158117
$ws.Cells.Item(94, 10).Value = 1128.6  # J94: was 1029.6364
$ws.Cells.Item(94, 11).Value = 20833616  # K94: was 13158117
$ws.Cells.Item(94, 12).Value = 1128.6  # L94: was 1029.6364
$ws.Cells.Item(94, 13).Value = -20833165  # M94: was -13157666
$ws.Cells.Item(94, 14).Value = -2030.6  # N94: was -1931.6364
$ws.Cells.Item(107, 8).Value = 1041.9032  # H107: was 1052.0625
$ws.Cells.Item(107, 9).Value = 803.1070999999999  # I107: was 846.38464
$ws.Cells.Item(107, 10).Value = 3270.6667  # J107: was 1943.3334
$ws.Cells.Item(107, 11).Value = 803.1070999999999  # K107: was 846.38464
$ws.Cells.Item(107, 12).Value = 3270.6667  # L107: was 1943.3334
$ws.Cells.Item(107, 13).Value = 1116.8929  # M107: was 1073.61536
$ws.Cells.Item(107, 14).Value = -7110.6667  # N107: was -5783.3334
$ws.Cells.Item(134, 8).Value = 4032.4736  # H134: was 4361.457
$ws.Cells.Item(134, 9).Value = 978.1142599999999  # I134: was 1051.5938
$ws.Cells.Item(134, 11).Value = 2934.34278  # K134: was 3154.7814
$ws.Cells.Item(134, 13).Value = -399.3427799999999  # M134: was -619.7814000000003
$ws.Cells.Item(136, 8).Value = 26691.334  # H136: was 25220
$ws.Cells.Item(136, 9).Value = 0  # I136: was 25000
$ws.Cells.Item(136, 10).Value = 26691.334  # J136: was 25330
$ws.Cells.Item(136, 11).Value = 0  # K136: was 25000
$ws.Cells.Item(136, 12).Value = 26691.334  # L136: was 25330
$ws.Cells.Item(136, 13).ClearContents()  # M136: remove (was -19900)
$ws.Cells.Item(136, 14).Value = -36891.334  # N136: was -35530

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1751.5957  # H31: was 1802.7333
$ws.Cells.Item(31, 9).Value = 1581.1464  # I31: was 1655.5526
$ws.Cells.Item(31, 10).Value = 2916.3333  # J31: was 2601.7144
$ws.Cells.Item(31, 11).Value = 1581.1464  # K31: was 1655.5526
$ws.Cells.Item(31, 12).Value = 2916.3333  # L31: was 2601.7144
$ws.Cells.Item(31, 13).Value = -1286.1464  # M31: was -1360.5526
$ws.Cells.Item(31, 14).Value = -3506.3333  # N31: was -3191.7144
$ws.Cells.Item(34, 8).Value = 1751.5957  # H34: was 1802.7333
$ws.Cells.Item(34, 9).Value = 1581.1464  # I34: was 1655.5526
$ws.Cells.Item(34, 10).Value = 2916.3333  # J34: was 2601.7144
$ws.Cells.Item(34, 11).Value = 1581.1464  # K34: was 1655.5526
$ws.Cells.Item(34, 12).Value = 2916.3333  # L34: was 2601.7144
$ws.Cells.Item(34, 13).Value = -1379.1464  # M34: was -1453.5526
$ws.Cells.Item(34, 14).Value = -3320.3333  # N34: was -3005.7144
$ws.Cells.Item(58, 8).Value = 881.32355  # H58: was 911.129
$ws.Cells.Item(58, 9).Value = 781.23334  # I58: was 804.3333
$ws.Cells.Item(58, 11).Value = 781.23334  # K58: was 804.3333
$ws.Cells.Item(58, 13).Value = -578.23334  # M58: was -601.3333
$ws.Cells.Item(132, 8).Value = 5420.6875  # H132: was 6654.16
$ws.Cells.Item(132, 9).Value = 6334.8335  # I132: was 8113.4443
$ws.Cells.Item(132, 10).Value = 2678.25  # J132: was 2901.7144
$ws.Cells.Item(132, 11).Value = 19004.5005  # K132: was 24340.3329
$ws.Cells.Item(132, 12).Value = 8034.75  # L132: was 8705.143199999999
$ws.Cells.Item(132, 13).Value = -16474.5005  # M132: was -21810.3329
$ws.Cells.Item(132, 14).Value = -13094.75  # N132: was -13765.1432
$ws.Cells.Item(136, 8).Value = 881.32355  # H136: was 911.129
$ws.Cells.Item(136, 9).Value = 781.23334  # I136: was 804.3333
$ws.Cells.Item(136, 11).Value = 2343.70002  # K136: was 2412.9999
$ws.Cells.Item(136, 13).Value = 206.2999799999998  # M136: was 137.0001000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(117, 8).Value = 1288.2858  # H117: was 1233.8572
$ws.Cells.Item(117, 9).Value = 629  # I117: was 625.6667
$ws.Cells.Item(117, 10).Value = 1552  # J117: was 1690
$ws.Cells.Item(117, 11).Value = 1887  # K117: was 1877.0001
$ws.Cells.Item(117, 12).Value = 4656  # L117: was 5070
$ws.Cells.Item(117, 13).Value = 1555  # M117: was 1564.9999
$ws.Cells.Item(117, 14).Value = -11540  # N117: was -11954
$ws.Cells.Item(129, 8).Value = 18940740  # H129: was 16027130
$ws.Cells.Item(129, 10).Value = 5557291  # J129: was 4387812.5
$ws.Cells.Item(129, 12).Value = 16671873  # L129: was 13163437.5
$ws.Cells.Item(129, 14).Value = -16681873  # N129: was -13173437.5
$ws.Cells.Item(131, 8).Value = 27779200  # H131: was 25001350
$ws.Cells.Item(131, 9).Value = 111111690  # I131: was 100000560
$ws.Cells.Item(131, 10).Value = 1704.3334  # J131: was 1615
$ws.Cells.Item(131, 11).Value = 333335070  # K131: was 300001680
$ws.Cells.Item(131, 12).Value = 5113.0002  # L131: was 4845
$ws.Cells.Item(131, 13).Value = -333330030  # M131: was -299996640
$ws.Cells.Item(131, 14).Value = -15193.0002  # N131: was -14925

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(26, 8).Value = 24800  # H26: was 24900
$ws.Cells.Item(26, 10).Value = 24800  # J26: was 24900
$ws.Cells.Item(26, 12).Value = 24800  # L26: was 24900
$ws.Cells.Item(26, 14).Value = -25360  # N26: was -25460
$ws.Cells.Item(50, 8).Value = 24800  # H50: was 24900
$ws.Cells.Item(50, 10).Value = 24800  # J50: was 24900
$ws.Cells.Item(50, 12).Value = 24800  # L50: was 24900
$ws.Cells.Item(50, 14).Value = -25796  # N50: was -25896

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2045.2222  # H7: was 1990.7
$ws.Cells.Item(7, 9).Value = 1915.4286  # I7: was 1863.5
$ws.Cells.Item(7, 11).Value = 1915.4286  # K7: was 1863.5
$ws.Cells.Item(7, 13).Value = -1803.4286  # M7: was -1751.5
$ws.Cells.Item(40, 8).Value = 1301.4062  # H40: was 1324.8529
$ws.Cells.Item(40, 9).Value = 1341.0714  # I40: was 1365
$ws.Cells.Item(40, 11).Value = 1341.0714  # K40: was 1365
$ws.Cells.Item(40, 13).Value = -1205.0714  # M40: was -1229
$ws.Cells.Item(82, 8).Value = 2162.8948  # H82: was 2157.0715
$ws.Cells.Item(82, 9).Value = 2371.7144  # I82: was 2188.889
$ws.Cells.Item(82, 10).Value = 2041.0834  # J82: was 2099.8
$ws.Cells.Item(82, 11).Value = 2371.7144  # K82: was 2188.889
$ws.Cells.Item(82, 12).Value = 2041.0834  # L82: was 2099.8
$ws.Cells.Item(82, 13).Value = -2010.7144  # M82: was -1827.889
$ws.Cells.Item(82, 14).Value = -2763.0834  # N82: was -2821.8
$ws.Cells.Item(85, 8).Value = 2162.8948  # H85: was 2157.0715
$ws.Cells.Item(85, 9).Value = 2371.7144  # I85: was 2188.889
$ws.Cells.Item(85, 10).Value = 2041.0834  # J85: was 2099.8
$ws.Cells.Item(85, 11).Value = 2371.7144  # K85: was 2188.889
$ws.Cells.Item(85, 12).Value = 2041.0834  # L85: was 2099.8
$ws.Cells.Item(85, 13).Value = -1123.7144  # M85: was -940.8890000000001
$ws.Cells.Item(85, 14).Value = -4537.0834  # N85: was -4595.8
$ws.Cells.Item(125, 8).Value = 35000  # H125: was 34230.77
$ws.Cells.Item(125, 10).Value = 35000  # J125: was 34230.77
$ws.Cells.Item(125, 12).Value = 35000  # L125: was 34230.77
$ws.Cells.Item(125, 14).Value = -44840  # N125: was -44070.77
$ws.Cells.Item(126, 8).Value = 2045.2222  # H126: was 1990.7
$ws.Cells.Item(126, 9).Value = 1915.4286  # I126: was 1863.5
$ws.Cells.Item(126, 11).Value = 5746.2858  # K126: was 5590.5
$ws.Cells.Item(126, 13).Value = -3276.2858  # M126: was -3120.5
$ws.Cells.Item(132, 8).Value = 34412.934  # H132: was 35840.332
$ws.Cells.Item(132, 9).Value = 2124.3333  # I132: was 2433.6667
$ws.Cells.Item(132, 10).Value = 54805.74  # J132: was 50157.477
$ws.Cells.Item(132, 11).Value = 6372.999899999999  # K132: was 7301.000100000001
$ws.Cells.Item(132, 12).Value = 164417.22  # L132: was 150472.431
$ws.Cells.Item(132, 13).Value = -3842.999899999999  # M132: was -4771.000100000001
$ws.Cells.Item(132, 14).Value = -169477.22  # N132: was -155532.431
$ws.Cells.Item(136, 8).Value = 6095.8696  # H136: was 5315.8887
$ws.Cells.Item(136, 9).Value = 7183.6113  # I136: was 6780.4736
$ws.Cells.Item(136, 10).Value = 2180  # J136: was 1837.5
$ws.Cells.Item(136, 11).Value = 21550.8339  # K136: was 20341.4208
$ws.Cells.Item(136, 12).Value = 6540  # L136: was 5512.5
$ws.Cells.Item(136, 13).Value = -19000.8339  # M136: was -17791.4208
$ws.Cells.Item(136, 14).Value = -11640  # N136: was -10612.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 141503.17  # H14: was 211254.75
$ws.Cells.Item(14, 9).Value = 177252.25  # I14: was 352504.5
$ws.Cells.Item(14, 11).Value = 177252.25  # K14: was 352504.5
$ws.Cells.Item(14, 13).Value = -177084.25  # M14: was -352336.5
$ws.Cells.Item(126, 8).Value = 47620764  # H126: was 45456228
$ws.Cells.Item(126, 9).Value = 58825016  # I126: was 55557012
$ws.Cells.Item(126, 11).Value = 176475048  # K126: was 166671036
$ws.Cells.Item(126, 13).Value = -176472578  # M126: was -166668566
$ws.Cells.Item(136, 8).Value = 592.8261  # H136: was 759.8
$ws.Cells.Item(136, 9).Value = 556  # I136: was 801.3333
$ws.Cells.Item(136, 10).Value = 661.875  # J136: was 697.5
$ws.Cells.Item(136, 11).Value = 1668  # K136: was 2403.9999
$ws.Cells.Item(136, 12).Value = 1985.625  # L136: was 2092.5
$ws.Cells.Item(136, 13).Value = 882  # M136: was 146.0001000000002
$ws.Cells.Item(136, 14).Value = -7085.625  # N136: was -7192.5
